# Saldo_guide.xlsx update script
# - Advance "Dt. Referencia" (column G) by one day for every data row (rows 2..310):
#     45386 (2024-04-04) -> 45387 (2024-04-05)
# - Adjust "Saldo Previsto" (D) / "Vl. Total" (H) for a handful of rows
# - Adjust "Vl. Projetado" (E) / "Vl. Total" (H) for row 120
# - Move the active selection to Q14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 310

# 1) Bump every reference date in column G from 45386 to 45387
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq 45386) {
        $cell.Value = 45387
    }
}

# 2) Row-specific balance corrections (Saldo Previsto / Vl. Total pairs)
$balanceChanges = @{
    6   = 35459.74
    12  = 67801.91
    23  = 1008.2
    47  = 3983.21
    55  = 24014.33
    121 = 1030.06
    126 = 21850.74
    129 = 84444.1
    130 = 16444.099999999999
    167 = 4425.8100000000004
    169 = 50.45
    187 = 244.23
    192 = 69934.92
}

foreach ($row in $balanceChanges.Keys) {
    $newValue = $balanceChanges[$row]
    $ws.Cells.Item($row, 4).Value = $newValue   # D - Saldo Previsto
    $ws.Cells.Item($row, 8).Value = $newValue   # H - Vl. Total
}

# 3) Row 120 special case: Vl. Projetado (E) changes, which feeds Vl. Total (H)
$ws.Cells.Item(120, 5).Value = 12945.98
$ws.Cells.Item(120, 8).Value = 109903.95

# 4) Update the visible selection
$ws.Range("Q14").Select()
